# Applies refreshed Yahoo Finance quote values (now behind a rate limiter)
# to the "Stock log", "Portfolio Summary" and "Total Return" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: Stock log (row 4 = SCHD) ---
$wsLog = $wb.Worksheets.Item("Stock log")
$wsLog.Range("F4").Value = 76.3
$wsLog.Range("J4").Value = 75002.89999999999
$wsLog.Range("K4").Value = 0.86
$wsLog.Range("M4").Value = 2.76

# --- Sheet: Portfolio Summary (row 3 = SCHD) ---
$wsSummary = $wb.Worksheets.Item("Portfolio Summary")
$wsSummary.Range("D3").Value = 75002.89999999999
$wsSummary.Range("G3").Value = 0.86
$wsSummary.Range("H3").Value = 2.76

# --- Sheet: Total Return ---
$wsReturn = $wb.Worksheets.Item("Total Return")
$wsReturn.Range("C2").Value = -1.73
$wsReturn.Range("D2").Value = 0.63
$wsReturn.Range("C3").Value = 0
$wsReturn.Range("D3").Value = -1.56
$wsReturn.Range("C4").Value = 0
$wsReturn.Range("D4").Value = 26.56
